# GUI: Updated the Statistics and made minor cosmetic adjustments to the Update suite.
#
# Adds a third CRUD "Update" source/column to the Folder Statistics sheet,
# mirroring the existing "Create" / "Read" rows, and nudges the saved
# selection to F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Register a brand-new external workbook link (the 3rd one) that will
#    back the new "Update" row, the same way the existing "Create"
#    (link #1) and "Read" (link #2) rows are backed by externalLink1.xml
#    / externalLink2.xml. We touch Sheet1/Sheet2/Sheet3 of the phantom
#    book once each (from a scratch cell well outside the used range) so
#    all three sheet names land in the exported externalLink3.xml, then
#    repoint the freshly created link at the real sibling workbook path
#    and clear the scratch cells again.
# ---------------------------------------------------------------------
$ws.Range("Z1000").Formula = "=[_UpdateLinkProbe.xlsx]Sheet1!`$A`$1"
$ws.Range("Z1001").Formula = "=[_UpdateLinkProbe.xlsx]Sheet2!`$A`$1"
$ws.Range("Z1002").Formula = "=[_UpdateLinkProbe.xlsx]Sheet3!`$A`$1"
$wb.ChangeLink("_UpdateLinkProbe.xlsx", "Update/_Test_Suite_Statistics_for_Folders.xlsx", 1) | Out-Null
$ws.Range("Z1000:Z1002").Clear() | Out-Null

# ---------------------------------------------------------------------
# 2. Add the new "Update" row label in A4, copying the look of the
#    "Create"/"Read" labels above it (A3) so the cell style matches.
# ---------------------------------------------------------------------
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A4").Value2 = "Update"

# ---------------------------------------------------------------------
# 3. Fill in the new row's per-suite-status figures (Ready to write,
#    Automated, Total, Automated cases, Total cases) for the Update
#    folder statistics, matching the values pulled from the new link.
# ---------------------------------------------------------------------
$ws.Range("B4").Value2 = 0
$ws.Range("C4").Value2 = 0
$ws.Range("D4").Value2 = 1
$ws.Range("E4").Value2 = 0
$ws.Range("F4").Value2 = 1

# ---------------------------------------------------------------------
# 4. Cosmetic: move the saved selection to F5.
# ---------------------------------------------------------------------
$ws.Range("F5").Select() | Out-Null
